# Apply cryptos list update (price/volume refresh + a few row re-orderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "34.645.24"
$ws.Cells.Item(2, 5).Value = "  +0.35%  "

$ws.Cells.Item(3, 4).Value = "1.814.35"
$ws.Cells.Item(3, 5).Value = "  +0.20%  "

$ws.Cells.Item(4, 5).Value = "  -0.07%  "

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "226.37"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.92%  "

$ws.Cells.Item(7, 5).Value = "  -0.07%  "

$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "38.25"
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +9.55%  "

$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.293"
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -3.58%  "

$ws.Cells.Item(10, 5).Value = "  -2.44%  "

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0971"
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +1.46%  "

$ws.Cells.Item(12, 4).Value = "2.074.44"
$ws.Cells.Item(12, 5).Value = "  +0.07%  "

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.32"
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +0.95%  "

$ws.Cells.Item(14, 4).Value = "1.819.55"
$ws.Cells.Item(14, 5).Value = "  +0.19%  "

$ws.Cells.Item(15, 5).Value = "  -2.34%  "

$ws.Cells.Item(16, 4).Value = "34.604.00"
$ws.Cells.Item(16, 5).Value = "  +0.25%  "

$ws.Cells.Item(17, 5).Value = "  -1.89%  "

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "68.90"
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -0.58%  "

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "244.88"
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.65%  "

$ws.Cells.Item(20, 5).Value = "  -2.60%  "

$ws.Cells.Item(21, 5).Value = "  -1.49%  "

$ws.Cells.Item(22, 5).Value = "  -0.07%  "

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.15"
$cell.Style = "Normal"

$ws.Cells.Item(24, 5).Value = "  +4.78%  "

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "172.08"
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.43%  "

$ws.Cells.Item(26, 5).Value = "  -2.45%  "

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "17.48"
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +3.91%  "

$ws.Cells.Item(28, 5).Value = "  +1.69%  "

$ws.Cells.Item(29, 5).Value = "  -0.06%  "

$ws.Cells.Item(30, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.93"
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -2.39%  "

$ws.Cells.Item(31, 2).Value = "Filecoin"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.82"
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -1.42%  "

$ws.Cells.Item(32, 5).Value = "  -2.33%  "

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.24"
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.97%  "

$ws.Cells.Item(34, 5).Value = "  -0.92%  "

$ws.Cells.Item(35, 4).Value = "1.367.10"
$ws.Cells.Item(35, 5).Value = "  -2.16%  "

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.657"
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -3.55%  "

$ws.Cells.Item(37, 5).Value = "  -0.75%  "

$ws.Cells.Item(38, 5).Value = "  -1.17%  "

$ws.Cells.Item(39, 5).Value = "  -5.29%  "

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.23"
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +7.65%  "

$ws.Cells.Item(41, 5).Value = "  +1.51%  "

$ws.Cells.Item(42, 2).Value = "Aave"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "81.21"
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -3.16%  "

$ws.Cells.Item(43, 2).Value = "ARBITRUM"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.943"
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -2.34%  "

$ws.Cells.Item(44, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "14.18"
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +6.28%  "

$ws.Cells.Item(45, 2).Value = "MXToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.78"
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -1.65%  "

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0502"
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -1.98%  "

$ws.Cells.Item(47, 4).Value = "1.975.16"
$ws.Cells.Item(47, 5).Value = "  +0.11%  "

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.79"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -3.74%  "

$ws.Cells.Item(49, 5).Value = "  -0.07%  "

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "103.04"
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -2.19%  "

$ws.Cells.Item(51, 2).Value = "BitcoinSV"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "49.20"
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -2.19%  "
